$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.154.71"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.22"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.39"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.96"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.17"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.650.79"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.11"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.44"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.143.51"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.49"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.43"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.51"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.41"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.118"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.266.98"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0177"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.542"
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.840"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.807"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.39"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.785.78"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.71"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +17.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0511"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.68"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0974"
$ws.Range("E51").Value = "  -1.11%  "
